# Update FRED WALCL data: append new weekly observations and refresh
# SeriesInfo metadata (custom fisher-index / BEA submodule data refresh).

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("Data")
$infoSheet = $wb.Worksheets.Item("SeriesInfo")

# --- Append new rows to the "Data" sheet (WALCL weekly series) ---
$newRows = @(
    @(104, 45189, 8024.09),
    @(105, 45196, 8002.064),
    @(106, 45203, 7955.782),
    @(107, 45210, 7952.054),
    @(108, 45217, 7933.162),
    @(109, 45224, 7907.83)
)

$lastStyledCell = $dataSheet.Cells.Item(103, 1)

foreach ($row in $newRows) {
    $r = $row[0]
    $dateVal = $row[1]
    $walclVal = $row[2]

    $dateCell = $dataSheet.Cells.Item($r, 1)
    $valCell = $dataSheet.Cells.Item($r, 2)

    # Copy the date column's number formatting/style from the prior row
    # before writing the new value so the appended cells keep the same
    # look (YYYY-MM-DD display, centered, bordered) as the rest of column A.
    $lastStyledCell.Copy() | Out-Null
    $dateCell.PasteSpecial(-4122) | Out-Null

    $dateCell.Value = $dateVal
    $valCell.Value = $walclVal
}

# --- Refresh SeriesInfo metadata to match the newly pulled data ---
# Force text format first so these date-shaped strings are stored as text
# (matching the source file's inlineStr cells) instead of being
# auto-converted to date serial numbers by Excel's input parser.
$realtimeStart = $infoSheet.Cells.Item(3, 2)
$realtimeStart.NumberFormat = "@"
$realtimeStart.Value = "2023-10-27"

$realtimeEnd = $infoSheet.Cells.Item(4, 2)
$realtimeEnd.NumberFormat = "@"
$realtimeEnd.Value = "2023-10-27"

$observationEnd = $infoSheet.Cells.Item(7, 2)
$observationEnd.NumberFormat = "@"
$observationEnd.Value = "2023-10-25"

$lastUpdated = $infoSheet.Cells.Item(14, 2)
$lastUpdated.NumberFormat = "@"
$lastUpdated.Value = "2023-10-26 15:33:02-05"

$infoSheet.Cells.Item(15, 2).Value = 93
